$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the current row 33, shifting the existing
# rows 33-89 down to rows 36-92 (a new week of "Maracuyá" price data is
# being prepended to this block).
$ws.Rows("33:35").Insert()

# Common (static) values shared by every data row in this block.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$codreg    = 15
$tipo      = "Fruta"
$productoId = 100108
$producto  = "Tropicales y subtropicales"
$categoriaId = 100108003
$categoria = "Maracuyá"
$variedad  = "Sin especificar"
$unidad    = "`$/caja 20 kilos"
$origen    = "Región de Arica y Parinacota"
$kgUnidad  = 20

# New rows of data (D=Fecha serial, L=Calidad, M=Volumen, N=Precio minimo,
# O=Precio maximo, P=Precio promedio ponderado, S=Precio $/Kg).
$newRows = @(
    @{ Row = 33; Fecha = 44624; Calidad = "Especial"; Volumen = 120; Min = 19000; Max = 20000; Prom = 19500; Kg = 975 },
    @{ Row = 34; Fecha = 44624; Calidad = "Primera";  Volumen = 120; Min = 16000; Max = 17000; Prom = 16500; Kg = 825 },
    @{ Row = 35; Fecha = 44624; Calidad = "Segunda";  Volumen = 130; Min = 14000; Max = 15000; Prom = 14500; Kg = 725 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.Min
    $ws.Cells.Item($row, 15).Value2 = $r.Max
    $ws.Cells.Item($row, 16).Value2 = $r.Prom
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.Kg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
